# Remodel of the datasheet: give ballistic horizontal-movement calculations
# their own "side Area" column (H) on the "quad" sheet, instead of reusing
# the top-area column (G). Mirrors how column E (drag-ish coefficient) is
# already computed from a constant-looking expression per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quad")
$ws.Activate()

# Row 1: standalone formula.
$ws.Range("H1").Formula = "=(0.058*0.069)"

# Rows 2-19: fill as one range write so Excel collapses it into a single
# shared formula group, same pattern already used for column E.
$ws.Range("H2:H19").Formula = "=(0.058*0.069)"

# Restore the cursor to where the author left it when the file was saved.
$ws.Range("J7").Select()
